# Apply the price / 1h-volume refresh captured in the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.758.08"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "1.642.61"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("E6").Value = "  -0.99%  "

$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0628"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("E9").Value = "  -1.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0840"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.57%  "

$ws.Range("D12").Value = "1.867.05"
$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("D13").Value = "1.644.59"
$ws.Range("E13").Value = "  +0.05%  "

$ws.Range("E14").Value = "  -1.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.527"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.24%  "

$ws.Range("D17").Value = "26.756.01"
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").Value = "0.0₃0737"
$ws.Range("E18").Value = "  -2.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.76%  "

$ws.Range("E20").Value = "  +0.31%  "

$ws.Range("E21").Value = "  -0.78%  "

$ws.Range("B22").Value = "Toncoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.28%  "

$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.65%  "

$ws.Range("E26").Value = "  +0.41%  "

$ws.Range("E27").Value = "  -2.19%  "

$ws.Range("E28").Value = "  -1.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0509"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.59%  "

$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("E32").Value = "  -2.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.24%  "

$ws.Range("D34").Value = "1.291.37"
$ws.Range("E34").Value = "  +0.65%  "

$ws.Range("E35").Value = "  -1.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.26%  "

$ws.Range("E37").Value = "  -6.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.535"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.825"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.53%  "

$ws.Range("E40").Value = "  +0.44%  "

$ws.Range("E41").Value = "  +0.10%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.15%  "

$ws.Range("D44").Value = "1.793.77"
$ws.Range("E44").Value = "  +0.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.30%  "

$ws.Range("E47").Value = "  -1.45%  "

$ws.Range("E48").Value = "  -1.63%  "

$ws.Range("E49").Value = "  +0.12%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0976"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.66%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.03%  "

